$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the two new rows (16, 17) by copying the formatting of row 15's
# A-column cell (style "1": bold/centered/bordered) before writing values.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

# Rows 8-15 are relabeled/renumbered (line7, line8 inserted ahead of the
# "extr" block, shifting it down by two), and rows 16-17 are brand new,
# reusing the extr7 / extr8 labels. Values below reflect the final state
# for each row after the edit.
$data = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $false },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
}
